$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Location - Innsbruck is a city in wich european country? / Austria
$ws.Range("A2").Value = "Innsbruck is a city in wich european country?"
$ws.Range("B2").Value = "Austria"
$ws.Range("C2").Value = "Location"

# Row 3: Location - What is the capital of Tyrol? / Innsbruck
$ws.Range("A3").Value = "What is the capital of Tyrol?"
$ws.Range("B3").Value = "Innsbruck"
$ws.Range("C3").Value = "Location"

# Row 4: Year - first F1 title / 1994
$ws.Range("A4").Value = "When did Michael Schumacher win his first F1 World Drivers Title?"
$ws.Range("B4").Value = 1994
$ws.Range("C4").Value = "Year"

# Row 5: Year - last F1 title / 2004
$ws.Range("A5").Value = "When did Michael Schumacher win his last F1 World Drivers Title?"
$ws.Range("B5").Value = 2004
$ws.Range("C5").Value = "Year"

# Row 6: Person - Grand Tour hosts / Richard Hammond
$ws.Range("A6").Value = "Who hosts The Grand Tour?"
$ws.Range("B6").Value = "Richard Hammond"
$ws.Range("C6").Value = "Person"

# Row 7: Person - 2022 F1 champion / Max Verstappen
$ws.Range("A7").Value = "Who was 2022 F1 World Drivers champion?"
$ws.Range("B7").Value = "Max Verstappen"
$ws.Range("C7").Value = "Person"

# Row 8: Person - Grand Tour hosts / Jeremy Clarkson
$ws.Range("A8").Value = "Who hosts The Grand Tour?"
$ws.Range("B8").Value = "Jeremy Clarkson"
$ws.Range("C8").Value = "Person"

# Column widths (bestFit-like autofit after content change)
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Selection update
$ws.Range("L13").Select() | Out-Null

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
